$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2:7 down to 3:8)
$ws.Rows("2:2").Insert()

# Grow the table (Table1) so it covers the newly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F8"))

# Drop the old "Need formatting" note that was on the algorithms row
# (shifted down, now row 3) ...
$ws.Range("F3").ClearContents()
# ... its text is folded into the "logic for computer science" row's note
# (shifted down, now row 4), replacing "Complete gedels proof"
$ws.Range("F4").Value = "Add algorithms, Complete gedels proof"

# Fill in the new entry ("ode", last updated 26/05/2024) on row 2
$ws.Range("A2").Value = "ode"
$ws.Range("B2").Value = 45438

# Copy the date style (numFmtId 14) from an existing "Last updated" cell
# onto B2 and C2 (C2 stays empty but picks up the same style)
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
